$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-29 Tuesday" "2025-07-30 Wednesday"

Replace-Text "46×92=" "86×91="
Replace-Text "24×14=" "40×17="
Replace-Text "34×97=" "39×57="
Replace-Text "38×85=" "15×76="
Replace-Text "65×84=" "21×48="
Replace-Text "70×80=" "18×43="
Replace-Text "46×65=" "51×79="
Replace-Text "62×76=" "25×62="
Replace-Text "91×61=" "44×74="
Replace-Text "56×54=" "47×85="
Replace-Text "54×31=" "92×26="
Replace-Text "52×93=" "19×49="
Replace-Text "76×71=" "12×61="
Replace-Text "16×57=" "71×25="
Replace-Text "40×59=" "66×16="
Replace-Text "63×62=" "40×86="
Replace-Text "25×14=" "74×75="
Replace-Text "16×58=" "48×37="
Replace-Text "75×33=" "15×17="
Replace-Text "34×64=" "42×19="
Replace-Text "90×34=" "20×37="
Replace-Text "83×49=" "94×15="
Replace-Text "93×41=" "15×99="
Replace-Text "46×33=" "43×33="
Replace-Text "20×62=" "82×87="
